$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Metadata" sheet updates
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version: 0.3.0 -> 0.4.0-snapshot-1
$meta.Range("B3").Value = "0.4.0-snapshot-1"

# Status: active -> draft
$meta.Range("B6").Value = "draft"

# Date: 2024-03-13T09:33:00+00:00 -> 2024-05-23T12:16:26+00:00
$meta.Range("B8").Value = "2024-05-23T12:16:26+00:00"

# Contact: "No display for ContactDetail" -> "ANS (https://esante.gouv.fr)"
$meta.Range("B10").Value = "ANS (https://esante.gouv.fr)"

# ---------------------------------------------------------------------------
# 2. "Elements" sheet: swap the two mapping columns AK ("Mapping: RIM Mapping")
#    and AL ("Mapping: Spécification métier vers l'extension ROR
#    LocationEquipmentLimit") so the business-mapping column comes first.
# ---------------------------------------------------------------------------
$els = $wb.Worksheets.Item("Elements")

# The "Spécification métier" column (formerly AL) is now first (AK) and is the
# wide one; the "RIM Mapping" column (formerly AK) is now second (AL) and is
# the narrow one - swap the column widths to match.
$els.Range("AK1").EntireColumn.ColumnWidth = 80.21875
$els.Range("AL1").EntireColumn.ColumnWidth = 24.98046875

for ($r = 1; $r -le 16; $r++) {
    $akCell = $els.Cells.Item($r, 37)  # column AK
    $alCell = $els.Cells.Item($r, 38)  # column AL

    $akVal = $akCell.Value2
    $alVal = $alCell.Value2

    $akCell.Value = $alVal
    $alCell.Value = $akVal
}
